$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing value in B2
$ws.Range("B2").Value = 37

# A3 gets a new "rank/count" value; also add new B3 value
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 30

# New rows 4 and 5 continue the count/value pairs, copying the style
# already used by A2:A3 (bold, bordered, centered) onto the new A cells
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A5").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A6").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 14

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 11

# The former row 3 content (A3=1, B3=3) moves down to row 6 with a new
# rescaled B value
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 7
